# Sync automático del tracker (cada 3h)
# Adds the newest tracked result as a new row at the bottom of the sheet
# (row 33), extending the used range from A1:H32 to A1:H33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append after the last used row (currently row 32 -> new data lands on
# row 33), rather than hard-coding the row number, so the sync keeps
# working as more rows accumulate over time.
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
$newRow = $lastRow + 1

# event_id / fecha are stored as plain text in the source feed (not as a
# number / date serial), so force text formatting before writing them —
# otherwise Excel would auto-coerce "14655144" to a number and
# "2025-09-20" to a date serial.
$idCell = $ws.Range("A" + $newRow)
$dateCell = $ws.Range("B" + $newRow)
$textRange = $ws.Range("A" + $newRow + ":B" + $newRow)
$textRange.NumberFormat = "@"

$idCell.Value = "14655144"
$dateCell.Value = "2025-09-20"

# Restore the default (unstyled) look now that the values are locked in
# as text, matching the rest of the sheet's unstyled data rows.
$textRange.Style = "Normal"

$ws.Range("C" + $newRow).Value = "Martin Damm Jr"
$ws.Range("D" + $newRow).Value = "Samir Banerjee"
$ws.Range("E" + $newRow).Value = "Gana Martin Damm Jr"
$ws.Range("F" + $newRow).Value = 2

# resultado / profit are still unknown (match not played yet) -- touch
# the cells so they exist in the sheet (matching the trailing empty
# columns used throughout the rest of the tracker) without altering
# their formatting.
$resultCell = $ws.Range("G" + $newRow)
$profitCell = $ws.Range("H" + $newRow)
$resultCell.Font.Bold = $false
$profitCell.Font.Bold = $false
